# Auto-generated edit script: updates the cryptos price/volume table
# to match the upstream GitHub Actions data refresh commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value2 = '63.130.05'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value2 = '  -0.57%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value2 = '2.611.48'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value2 = '  -2.16%  '
$ws.Range('E4').Value2 = '  +0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value2 = '605.86'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value2 = '  +1.93%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value2 = '145.46'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value2 = '  +0.91%  '
$ws.Range('E7').Value2 = '  -0.01%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value2 = '0.584'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value2 = '  -0.51%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value2 = '2.611.36'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value2 = '  -2.12%  '
$ws.Range('E10').Value2 = '  +0.74%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value2 = '5.50'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value2 = '  -3.46%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value2 = '0.372'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value2 = '  +4.33%  '
$ws.Range('E13').Value2 = '  -0.60%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value2 = '27.13'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value2 = '  -1.49%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value2 = '3.080.02'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value2 = '  -2.09%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value2 = '62.986.98'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value2 = '  -0.64%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value2 = '0.0000146'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value2 = '  +0.96%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value2 = '2.588.38'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value2 = '  -2.85%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value2 = '11.48'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value2 = '  -0.68%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value2 = '4.50'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value2 = '  +1.63%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value2 = '341.75'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value2 = '  +0.76%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value2 = '6.86'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value2 = '  +0.30%  '
$ws.Range('E23').Value2 = '  -0.09%  '
$ws.Range('E24').Value2 = '  -1.38%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value2 = '66.10'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value2 = '  -2.08%  '
$ws.Range('E26').Value2 = '  -0.39%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value2 = '1.59'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value2 = '  +3.13%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value2 = '9.00'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value2 = '  +5.51%  '
$ws.Range('B29').Value2 = 'Kaspa'
$ws.Range('C29').Value2 = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value2 = '0.162'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value2 = '  -2.46%  '
$ws.Range('B30').Value2 = 'Bittensor'
$ws.Range('C30').Value2 = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value2 = '542.69'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value2 = '  +0.29%  '
$ws.Range('E31').Value2 = '  +0.15%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value2 = '7.80'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value2 = '  -0.88%  '
$ws.Range('E33').Value2 = '  +1.81%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value2 = '0.0₃0840'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value2 = '  +3.24%  '
$ws.Range('E35').Value2 = '  -5.23%  '
$ws.Range('E36').Value2 = '  +1.08%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value2 = '168.09'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value2 = '  -2.82%  '
$ws.Range('E38').Value2 = '  +0.08%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value2 = '0.401'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value2 = '  -1.31%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value2 = '1.93'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value2 = '  +5.25%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value2 = '18.90'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value2 = '  -1.57%  '
$ws.Range('E42').Value2 = '  +0.08%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value2 = '164.95'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value2 = '  -5.62%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value2 = '39.66'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value2 = '  -1.34%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value2 = '3.75'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value2 = '  -0.55%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value2 = '21.73'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value2 = '  -2.78%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value2 = '0.0562'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value2 = '  -0.27%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value2 = '0.623'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value2 = '  -1.86%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value2 = '0.0242'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value2 = '  +0.62%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value2 = '0.0956'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value2 = '  -0.83%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value2 = '1.92'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value2 = '  +10.75%  '

Write-Output "Updated $($wb.Name): 91 cells (cryptos list refresh)"
